$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column D (Khoa) values: replace placeholder letters (A-I) with actual
# department / major names.
$ws.Range("D2").Value = "Kỹ thuật điện tử"
$ws.Range("D3").Value = "Khoa học máy tính"
$ws.Range("D4").Value = "Kỹ thuật viễn thông"
$ws.Range("D5").Value = "Kỹ thuật xây dựng"
$ws.Range("D6").Value = "Vật liệu và linh kiện nano"
$ws.Range("D7").Value = "An toàn thông tin"
$ws.Range("D8").Value = "Cơ kỹ thuật"
$ws.Range("D9").Value = "Kỹ thuật phần mềm"
$ws.Range("D10").Value = "Mạng máy tính và truyền thông dữ liệu"

# Resize column D to fit the new, longer text (bestFit width, like columns B/C).
$ws.Columns.Item(4).EntireColumn.AutoFit() | Out-Null

# Update the active selection, as last left by the editing user.
$ws.Range("D7").Select() | Out-Null
